$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name column (A) values to include links, matching the new source data
$ws.Range("A6").Value = "WordPress - https://en.wordpress.com/features/"
$ws.Range("A2").Value = "concrete5 - https://www.concrete5.org/"
$ws.Range("A5").Value = "Liferay Portal (ENTERPRISE SOLUTION ONLY) - https://www.liferay.com/product/features"

# Adjust the view: scroll so row 3 is the top-left visible row, select A6
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A6").Select()
